$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing scores on row 74 (sleep, 2025-02-25) ---
$ws.Range("F74").Value = $false

# --- Fix existing scores on row 75 (activity, 2025-02-25) ---
$ws.Range("C75").Value = $true
$ws.Range("F75").Value = $true
$ws.Range("I75").Value = $true
$ws.Range("M75").Value = $true

# --- Append new daily-score rows for 2025-02-26 ---

# Row 77: sleep
$ws.Range("A77").NumberFormat = "@"
$ws.Range("A77").Value = "2025-02-26"
$ws.Range("A77").ClearFormats()
$ws.Range("B77").Value = "sleep"
$ws.Range("C77").Value = $false
$ws.Range("D77").Value = $false
$ws.Range("E77").Value = $true
$ws.Range("F77").Value = $false
$ws.Range("G77").Value = $true
$ws.Range("H77").Value = $true
$ws.Range("I77").Value = $true
$ws.Range("J77").Value = $true
$ws.Range("K77").Value = $true
$ws.Range("L77").Value = $true
$ws.Range("M77").Value = $true
$ws.Range("N77").Value = $false
$ws.Range("O77").Value = $false

# Row 78: activity
$ws.Range("A78").NumberFormat = "@"
$ws.Range("A78").Value = "2025-02-26"
$ws.Range("A78").ClearFormats()
$ws.Range("B78").Value = "activity"
$ws.Range("C78").Value = $true
$ws.Range("D78").Value = $false
$ws.Range("E78").Value = $false
$ws.Range("F78").Value = $false
$ws.Range("G78").Value = $true
$ws.Range("H78").Value = $true
$ws.Range("I78").Value = $false
$ws.Range("J78").Value = $true
$ws.Range("K78").Value = $false
$ws.Range("L78").Value = $false
$ws.Range("M78").Value = $true
$ws.Range("N78").Value = $false
$ws.Range("O78").Value = $false

# Row 79: weekly_activity
$ws.Range("A79").NumberFormat = "@"
$ws.Range("A79").Value = "2025-02-26"
$ws.Range("A79").ClearFormats()
$ws.Range("B79").Value = "weekly_activity"
$ws.Range("C79").Value = $true
$ws.Range("D79").Value = $false
$ws.Range("E79").Value = $true
$ws.Range("F79").Value = $true
$ws.Range("G79").Value = $true
$ws.Range("H79").Value = $true
$ws.Range("I79").Value = $true
$ws.Range("J79").Value = $true
$ws.Range("K79").Value = $false
$ws.Range("L79").Value = $false
$ws.Range("M79").Value = $true
$ws.Range("N79").Value = $true
$ws.Range("O79").Value = $false
